# Add the new "2022-Q3" quarter sheet right after "总计" (总计 stays first).
$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Copy the header/label formatting (bold + border style used by the other
# quarter sheets) from the "2022-Q2" sheet onto the new sheet, then fill in
# the 2022-Q3 fund data.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Range("A1:H3").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "519702"
$newSheet.Range("C2").Value = "交银趋势优先混合A"
$newSheet.Range("D2").Value = "83.94"
$newSheet.Range("E2").Value = "81.61"
$newSheet.Range("F2").Value = "2.34"
$newSheet.Range("G2").Value = "1.9642"
$newSheet.Range("H2").Value = 7

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "013430"
$newSheet.Range("C3").Value = "交银趋势优先混合C"
$newSheet.Range("D3").Value = "12.61"
$newSheet.Range("E3").Value = "81.61"
$newSheet.Range("F3").Value = "2.34"
$newSheet.Range("G3").Value = "0.2951"
$newSheet.Range("H3").Value = 7

# Update the "总计" summary sheet: shift each quarter's row down to make
# room for the new 2022-Q3 figures at the top of the data (row 2).
$totalSheet.Range("B4").Value = "2021-Q4"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 1.97

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 3
$totalSheet.Range("D5").Value = 1.97

$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)

$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 2.55

$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 4
$totalSheet.Range("D3").Value = 4.22

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 2.26
